$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- New rows 15-18 : data for "1_voice" / "1_background" / etc -----------
$ws.Range("A15").Value = 43369
$ws.Range("A15").NumberFormat = "d-mmm"
$ws.Range("B15").Value = 52
$ws.Range("C15").Value = "0_voice"
$ws.Range("D15").Value = [double]"5.8499455248008101E-3"
$ws.Range("E15").Value = -4.9414299214992896
$ws.Range("F15").Value = -2.9187624359262498
$ws.Range("G15").Value = 8.0849444575939504
$ws.Range("H15").Value = [double]"-6.4876311729494796E-4"
$ws.Range("I15").Value = "Unet trained on full chime and librispeech dataset (10 epochs)"

$ws.Range("C16").Value = "0_background"
$ws.Range("D16").Value = [double]"5.8499455248008101E-3"
$ws.Range("E16").Value = 0.603653156846512
$ws.Range("F16").Value = 4.9341611496209197
$ws.Range("G16").Value = 8.0871218884956892
$ws.Range("H16").Value = [double]"1.54551858325779E-2"

$ws.Range("C17").Value = "1_voice"
$ws.Range("D17").Value = [double]"3.9196163864368798E-3"
$ws.Range("E17").Value = 0.10210658794458299
$ws.Range("F17").Value = 6.2019812058338104
$ws.Range("G17").Value = 4.9962338577947696
$ws.Range("H17").Value = 5.0428877463265804

$ws.Range("C18").Value = "1_background"
$ws.Range("D18").Value = [double]"3.9196163864368798E-3"
$ws.Range("E18").Value = 4.17242882565266
$ws.Range("F18").Value = 13.924662842344301
$ws.Range("G18").Value = 5.65928442136721
$ws.Range("H18").Value = 3.58423085463872

# --- Row 14 gains A/B/I placeholder cells (to be merged upward) -----------
$ws.Range("A14").Value = 43368
$ws.Range("A14").NumberFormat = "d-mmm"
$ws.Range("B14").Value = 51
$ws.Range("I14").Value = ""

# --- Center-align the Experiment (B) / Date (A) / Description (I) columns -
$ws.Range("B9:B18").HorizontalAlignment = -4108
$ws.Range("A11:A18").HorizontalAlignment = -4108
$ws.Range("I13:I18").HorizontalAlignment = -4108

# --- Column C width (bestFit) ----------------------------------------------
$ws.Columns.Item(3).ColumnWidth = 13.42578125

# --- Merge repeated Date / Experiment / Description cells -----------------
$ws.Range("B9:B10").Merge()
$ws.Range("A11:A12").Merge()
$ws.Range("B11:B12").Merge()
$ws.Range("A13:A14").Merge()
$ws.Range("B13:B14").Merge()
$ws.Range("I13:I14").Merge()
$ws.Range("A15:A18").Merge()
$ws.Range("B15:B18").Merge()
$ws.Range("I15:I18").Merge()

# --- Selection -------------------------------------------------------------
$ws.Range("I26").Select()
